$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C14").Value = "Rev_03_03"
$ws.Range("C15").Value = "Rev_04_01"
$ws.Range("C16").Value = "Rev_05_01"
$ws.Range("C17").Value = "Rev_06_01"
$ws.Range("C18").Value = "Rev_07_01"
$ws.Range("C19").Value = "Rev_08_01"
$ws.Range("C20").Value = "Rev_09_01"
$ws.Range("C21").Value = "Rev_10_01"
$ws.Range("C22").Value = "Rev_11_01"
$ws.Range("C23").Value = "Rev_12_01"

$ws.Range("C24").Select()
